$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 19, pushing existing rows 19-86 down to 20-87.
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with the new record.
$ws.Range("A19").Value = 10
$ws.Range("B19").Value = "Vega Modelo de Temuco"
$ws.Range("C19").Value = "La Araucanía"
$ws.Range("D19").Value = 44497
$ws.Range("E19").Value = 9
$ws.Range("F19").Value = "Fruta"
$ws.Range("G19").Value = 100107
$ws.Range("H19").Value = "Otros"
$ws.Range("I19").Value = 100107002
$ws.Range("J19").Value = "Chirimoya"
$ws.Range("K19").Value = "Cultivar IV Región"
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 140
$ws.Range("N19").Value = 3000
$ws.Range("O19").Value = 3000
$ws.Range("P19").Value = 3000
$ws.Range("Q19").Value = "`$/kilo (en caja de 15 kilos)"
$ws.Range("R19").Value = "Provincia del Elquí"
$ws.Range("S19").Value = 3000
$ws.Range("T19").Value = 1
